$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 8 (year 2025) metrics with new recalculated figures
$ws.Range("C8").Value = 1133
$ws.Range("D8").Value = 187
$ws.Range("E8").Value = 946
$ws.Range("F8").Value = 7.670221493027071
$ws.Range("G8").Value = 83.49514563106796
$ws.Range("H8").Value = 16.50485436893204
